# Apply changes described by the diff:
# - B5 changes from NET1014 to ENG1044
# - New rows 6-13 added for Semester 2 and Semester 3 modules
# - Selection moves to D16
# - Dimension grows to A1:B13 (handled automatically by Excel when cells are written)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Semester 2 modules (entered first so shared strings register in this order)
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "SEG1201"

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "OSS1014"

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "WEB1201"

$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "PRG1203"

# Update B5 (was NET1014) to ENG1044
$ws.Range("B5").Value = "ENG1044"

# Semester 3 modules
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "SEG1201"

$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "OSS1014"

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "WEB1201"

$ws.Range("A13").Value = 3
$ws.Range("B13").Value = "NET1014"

# Move the active selection to D16, matching the final saved view state
$ws.Range("D16").Select()
